$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color helpers (Excel COM colors are BGR, not RGB)
$navy   = 0x5D3617   # RGB 17365D
$white  = 0xFFFFFF
$black  = 0x000000
$orange = 0x8FBFFA   # RGB FABF8F
$yellow = 0xCCFFFF   # RGB FFFFCC
$cyan   = 0xFFFFCC   # RGB CCFFFF
$lime   = 0x00F0F0   # RGB F0F000

# ---------------------------------------------------------------------------
# 1. Header text (A1:G1 / A2:G2 / A3:G3)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Delfiniti de México S.A. de C.V."
$ws.Range("A2").Value = "REPORTE DE RESERVACIONES"
$ws.Range("A3").Value = "Del 17-08-2022 al 17-08-2022"

$hdr1 = $ws.Range("A1:G1")
$hdr1.Interior.Color = $navy
$hdr1.Interior.PatternColor = $white
$hdr1.Font.Color = $white
$hdr1.Font.Size = 26
$hdr1.Font.Bold = $true
$hdr1.Font.Name = "Calibri"
$hdr1.HorizontalAlignment = -4108
$hdr1.VerticalAlignment = -4108
$hdr1.WrapText = $true

$hdr2 = $ws.Range("A2:G2")
$hdr2.Interior.Color = $navy
$hdr2.Interior.PatternColor = $white
$hdr2.Font.Color = $white
$hdr2.Font.Size = 20
$hdr2.Font.Bold = $false
$hdr2.Font.Name = "Calibri"
$hdr2.HorizontalAlignment = -4108
$hdr2.VerticalAlignment = -4108
$hdr2.WrapText = $true

$hdr3 = $ws.Range("A3:G3")
$hdr3.Interior.Color = $navy
$hdr3.Interior.PatternColor = $white
$hdr3.Font.Color = $white
$hdr3.Font.Size = 12
$hdr3.Font.Bold = $false
$hdr3.Font.Name = "Calibri"
$hdr3.HorizontalAlignment = -4108
$hdr3.VerticalAlignment = -4108
$hdr3.WrapText = $true

Write-Host "headers done"

# ---------------------------------------------------------------------------
# 2. Drop the old "PROGRAMA / PAGADOS / PENDIENTES / CORTESIAS / TOTAL" block
#    that used to sit at rows 8:9 right under the header - it gets rebuilt
#    further down (rows 18:21) with the new per-event breakdown above it.
# ---------------------------------------------------------------------------
$ws.Range("A8:O9").EntireRow.Delete() | Out-Null
Write-Host "old total rows removed"

# ---------------------------------------------------------------------------
# 3. Event block #1 - "FAMILY CLUB MED 10:00:00" (rows 5-8)
# ---------------------------------------------------------------------------
$ws.Range("A5:C5").Merge()
$ws.Range("A5").Value = "FAMILY CLUB MED 10:00:00"
$blk1 = $ws.Range("A5:C5")
$blk1.Interior.Color = $orange
$blk1.Interior.PatternColor = $black
$blk1.HorizontalAlignment = 1
$blk1.VerticalAlignment = -4108

$ws.Range("A6:B6").Merge()
$ws.Range("A6").Value = "CLIENTE"
$ws.Range("C6").Value = "ORIGEN"
$ws.Range("D6").Value = "PAX"
$ws.Range("E6").Value = "AGENTE/AGENCIA"
$ws.Range("F6").Value = "T. PAGO"
$row6 = $ws.Range("A6:F6")
$row6.Interior.Color = $yellow
$row6.Interior.PatternColor = $black
$row6.HorizontalAlignment = 1
$row6.VerticalAlignment = -4108

$ws.Range("A7:B7").Merge()
$ws.Range("A7").Value = "JORGE"
$ws.Range("D7").Value = 1
$ws.Range("F7").Value = "efectivo"
$row7 = $ws.Range("A7:F7")
$row7.Interior.Color = $cyan
$row7.Interior.PatternColor = $black
$row7.HorizontalAlignment = 1
$row7.VerticalAlignment = -4108

$ws.Range("D8").Formula = "=SUM(D7:D7)"

Write-Host "block1 done"

# ---------------------------------------------------------------------------
# 4. Event block #2 - "ENCUENTRO 10:00:00" (rows 10-13)
# ---------------------------------------------------------------------------
$ws.Range("A10:C10").Merge()
$ws.Range("A10").Value = "ENCUENTRO 10:00:00"
$blk2 = $ws.Range("A10:C10")
$blk2.Interior.Color = $orange
$blk2.Interior.PatternColor = $black
$blk2.HorizontalAlignment = 1
$blk2.VerticalAlignment = -4108

$ws.Range("A11:B11").Merge()
$ws.Range("A11").Value = "CLIENTE"
$ws.Range("C11").Value = "ORIGEN"
$ws.Range("D11").Value = "PAX"
$ws.Range("E11").Value = "AGENTE/AGENCIA"
$ws.Range("F11").Value = "T. PAGO"
$row11 = $ws.Range("A11:F11")
$row11.Interior.Color = $yellow
$row11.Interior.PatternColor = $black
$row11.HorizontalAlignment = 1
$row11.VerticalAlignment = -4108

$ws.Range("A12:B12").Merge()
$ws.Range("A12").Value = "TEST"
$ws.Range("D12").Value = 5
$row12 = $ws.Range("A12:F12")
$row12.Interior.Color = $cyan
$row12.Interior.PatternColor = $black
$row12.HorizontalAlignment = 1
$row12.VerticalAlignment = -4108

$ws.Range("D13").Formula = "=SUM(D12:D12)"

Write-Host "block2 done"

# ---------------------------------------------------------------------------
# 5. Summary table (rows 18-21) - PROGRAMA / PAGADOS / PENDIENTES / CORTESIAS / TOTAL
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "PROGRAMA"
$ws.Range("B18").Value = "PAGADOS"
$ws.Range("C18").Value = "PENDIENTES"
$ws.Range("D18").Value = "CORTESIAS"
$ws.Range("E18").Value = "TOTAL"
$row18 = $ws.Range("A18:E18")
$row18.Interior.Color = $lime
$row18.Interior.PatternColor = $black
$row18.HorizontalAlignment = 1
$row18.VerticalAlignment = -4108

$ws.Range("A19").Value = "FAMILY CLUB MED"
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 0
$ws.Range("E19").Formula = "=SUM(B19:D19)"

$ws.Range("A20").Value = "ENCUENTRO"
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 0
$ws.Range("E20").Formula = "=SUM(B20:D20)"

$ws.Range("B21").Formula = "=SUM(B19:B20)"
$ws.Range("C21").Formula = "=SUM(C19:C20)"
$ws.Range("D21").Formula = "=SUM(D19:D20)"
$ws.Range("E21").Formula = "=SUM(E19:E20)"
$totalRow = $ws.Range("A21:D21")
$totalRow.Interior.Color = $orange
$totalRow.Interior.PatternColor = $black
$totalRow.HorizontalAlignment = 1
$totalRow.VerticalAlignment = -4108

Write-Host "summary table done"

# ---------------------------------------------------------------------------
# 6. Column widths (values chosen so the saved OOXML <col width> lands on the
#    target figure - the ColumnWidth COM property only has 1/6-character
#    granularity, so these are the closest achievable inputs)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(3).ColumnWidth = 11.666666666666666
$ws.Columns.Item(4).ColumnWidth = 8.333333333333334
$ws.Columns.Item(5).ColumnWidth = 23.5
$ws.Columns.Item(6).ColumnWidth = 11
$ws.Columns.Item(7).ColumnWidth = 12

Write-Host "column widths done"

# ---------------------------------------------------------------------------
# 7. Selection - mirrors the new bottom-right total cell block
# ---------------------------------------------------------------------------
$ws.Range("A21:D21").Select() | Out-Null

Write-Host "selection done"

